$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised historical values in rows 419-421 (col B & D changed) ---
$ws.Range("B419").Value = 5245350000000
$ws.Range("D419").Value = 1051425192431.045

$ws.Range("B420").Value = 5322265000000
$ws.Range("D420").Value = 1052788107765.953

$ws.Range("B421").Value = 5421638000000
$ws.Range("D421").Value = 1131984131955.319

# --- Append 3 new monthly rows (422-424), matching col A's date style (s="2") ---
$ws.Range("A421").Copy()
$ws.Range("A422:A424").PasteSpecial(-4122)

$ws.Range("A422").Value = 45108
$ws.Range("B422").Value = 5501072000000
$ws.Range("C422").Value = 0.2116894938504202
$ws.Range("D422").Value = 1164519147314.719

$ws.Range("A423").Value = 45139
$ws.Range("B423").Value = 5591097000000
$ws.Range("C423").Value = 0.2018733850129199
$ws.Range("D423").Value = 1128693677325.582

$ws.Range("A424").Value = 45170
$ws.Range("B424").Value = 5656835000000
$ws.Range("C424").Value = 0.1987281399046105
$ws.Range("D424").Value = 1124172297297.297
